# Auto-generated Excel COM-interop script
# Applies cell value updates to multiple worksheets per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 999
$ws.Range("I2").Value = 999
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 999
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -886
$ws.Range("H18").Value = 910.36365
$ws.Range("I18").Value = 910.36365
$ws.Range("K18").Value = 910.36365
$ws.Range("M18").Value = -626.36365
$ws.Range("H40").Value = 4554.4165
$ws.Range("J40").Value = 4628.1113
$ws.Range("L40").Value = 4628.1113
$ws.Range("N40").Value = -4978.1113
$ws.Range("H64").Value = 9378.75
$ws.Range("H67").Value = 9378.75
$ws.Range("H95").Value = 32000
$ws.Range("J95").Value = 32000
$ws.Range("L95").Value = 32000
$ws.Range("N95").Value = -37492
$ws.Range("H113").Value = 2437.2727
$ws.Range("I113").Value = 1801.375
$ws.Range("K113").Value = 1801.375
$ws.Range("M113").Value = 1452.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 474.83334
$ws.Range("I4").Value = 474.83334
$ws.Range("K4").Value = 474.83334
$ws.Range("M4").Value = -358.83334
$ws.Range("H16").Value = 7143203
$ws.Range("I16").Value = 403.33334
$ws.Range("K16").Value = 403.33334
$ws.Range("M16").Value = -116.33334
$ws.Range("H26").Value = 4500
$ws.Range("I26").Value = 4500
$ws.Range("K26").Value = 4500
$ws.Range("M26").Value = -4170
$ws.Range("H31").Value = 6617.5
$ws.Range("I31").Value = 6617.5
$ws.Range("K31").Value = 6617.5
$ws.Range("M31").Value = -6323.5
$ws.Range("H39").Value = 2155.3333
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H96").Value = 13905.333
$ws.Range("J96").Value = 13905.333
$ws.Range("L96").Value = 13905.333
$ws.Range("N96").Value = -19397.333
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 945.5
$ws.Range("I7").Value = 140.25
$ws.Range("K7").Value = 140.25
$ws.Range("M7").Value = -27.25
$ws.Range("H48").Value = 199999
$ws.Range("J48").Value = 199999
$ws.Range("L48").Value = 199999
$ws.Range("N48").Value = -200829
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0
$ws.Range("H106").Value = 23900.166
$ws.Range("J106").Value = 23900.166
$ws.Range("L106").Value = 23900.166
$ws.Range("N106").Value = -26424.166
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 866.46155
$ws.Range("I2").Value = 698.75
$ws.Range("J2").Value = 1134.8
$ws.Range("K2").Value = 698.75
$ws.Range("L2").Value = 1134.8
$ws.Range("M2").Value = -585.75
$ws.Range("N2").Value = -1360.8
$ws.Range("H12").Value = 636
$ws.Range("I12").Value = 475.75
$ws.Range("K12").Value = 475.75
$ws.Range("M12").Value = -305.75
$ws.Range("H17").Value = 233.33333
$ws.Range("I17").Value = 233.33333
$ws.Range("K17").Value = 233.33333
$ws.Range("M17").Value = -59.33332999999999
$ws.Range("H19").Value = 1169.8572
$ws.Range("I19").Value = 70
$ws.Range("J19").Value = 1609.8
$ws.Range("K19").Value = 70
$ws.Range("L19").Value = 1609.8
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = -1949.8
$ws.Range("H24").Value = 1169.8572
$ws.Range("I24").Value = 70
$ws.Range("J24").Value = 1609.8
$ws.Range("K24").Value = 70
$ws.Range("L24").Value = 1609.8
$ws.Range("M24").Value = 100
$ws.Range("N24").Value = -1949.8
$ws.Range("H33").Value = 24969.8
$ws.Range("I33").Value = 1616.6666
$ws.Range("K33").Value = 1616.6666
$ws.Range("M33").Value = -1237.6666
$ws.Range("H36").Value = 36665.668
$ws.Range("I36").Value = 19998.5
$ws.Range("J36").Value = 70000
$ws.Range("K36").Value = 19998.5
$ws.Range("L36").Value = 70000
$ws.Range("M36").Value = -19610.5
$ws.Range("N36").Value = -70776
$ws.Range("H40").Value = 36665.668
$ws.Range("I40").Value = 19998.5
$ws.Range("J40").Value = 70000
$ws.Range("K40").Value = 19998.5
$ws.Range("L40").Value = 70000
$ws.Range("M40").Value = -19838.5
$ws.Range("N40").Value = -70320
$ws.Range("H51").Value = 19647.5
$ws.Range("I51").Value = 19647.5
$ws.Range("K51").Value = 19647.5
$ws.Range("M51").Value = -18911.5
$ws.Range("H61").Value = 19647.5
$ws.Range("I61").Value = 19647.5
$ws.Range("K61").Value = 19647.5
$ws.Range("M61").Value = -19299.5
$ws.Range("H96").Value = 37837.332
$ws.Range("J96").Value = 37837.332
$ws.Range("L96").Value = 37837.332
$ws.Range("N96").Value = -43329.332
$ws.Range("H105").Value = 3061.25
$ws.Range("I105").Value = 3081.6667
$ws.Range("K105").Value = 3081.6667
$ws.Range("M105").Value = -1334.6667
$ws.Range("H134").Value = 2641.4
$ws.Range("I134").Value = 2470.8462
$ws.Range("K134").Value = 7412.5386
$ws.Range("M134").Value = -4877.5386
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 423.4
$ws.Range("I8").Value = 423.4
$ws.Range("K8").Value = 1270.2
$ws.Range("M8").Value = -1131.2
$ws.Range("H114").Value = 1533.3125
$ws.Range("I114").Value = 812
$ws.Range("J114").Value = 1861.1818
$ws.Range("K114").Value = 2436
$ws.Range("L114").Value = 5583.5454
$ws.Range("M114").Value = 818
$ws.Range("N114").Value = -12091.5454
$ws.Range("H117").Value = 4516.364
$ws.Range("J117").Value = 7633.3335
$ws.Range("L117").Value = 22900.0005
$ws.Range("N117").Value = -29784.0005
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 51020.5
$ws.Range("J58").Value = 52000
$ws.Range("L58").Value = 52000
$ws.Range("N58").Value = -52554
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 412.5
$ws.Range("I9").Value = 216.66667
$ws.Range("K9").Value = 216.66667
$ws.Range("M9").Value = 7.333329999999989
$ws.Range("H22").Value = 677.55554
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 677.55554
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").Value = 677.55554
$ws.Range("N22").Value = -1267.55554
$ws.Range("H27").Value = 677.55554
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 677.55554
$ws.Range("K27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("M27").Value = 677.55554
$ws.Range("N27").Value = -891.55554
$ws.Range("H35").Value = 10068.214
$ws.Range("I35").Value = 3870
$ws.Range("K35").Value = 3870
$ws.Range("M35").Value = -3534
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 30999
$ws.Range("J95").Value = 30999
$ws.Range("L95").Value = 30999
$ws.Range("N95").Value = -36491
$ws.Range("H136").Value = 1936.1428
$ws.Range("I136").Value = 1924.75
$ws.Range("K136").Value = 5774.25
$ws.Range("M136").Value = -3224.25
